# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.192.35"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -1.08%  "

# Row 3
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.578.90"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -1.62%  "

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.10"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.25%  "

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.07"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -3.77%  "

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.573.49"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -1.60%  "

# Row 8
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -3.31%  "

# Row 9
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("E10").Value = "  +1.51%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.654"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.90%  "

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.33"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -6.03%  "

# Row 13
$ws.Range("E13").Value = "  -2.53%  "

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.61"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -3.09%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.151.38"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -1.76%  "

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.71"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -3.56%  "

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.579.99"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -1.66%  "

# Row 18
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.039.51"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -1.25%  "

# Row 19
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -1.94%  "

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.121"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -1.43%  "

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.19%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "499.48"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +1.88%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.48"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +2.97%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.07"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -1.63%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.68"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +6.35%  "

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.40"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -2.56%  "

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.71"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +2.26%  "

# Row 28
$ws.Range("E28").Value = "  -5.25%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.37"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -2.15%  "

# Row 30
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -3.36%  "

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.72"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -3.37%  "

# Row 32
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.83"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +4.40%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.89"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -2.60%  "

# Row 34
$ws.Range("E34").Value = "  -4.09%  "

# Row 35
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "569.84"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -7.27%  "

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.31"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +12.21%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.08"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -4.03%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.411"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +0.33%  "

# Row 39
$ws.Range("E39").Value = "  -0.05%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0790"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -7.73%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.20"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -2.44%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.46"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -3.37%  "

# Row 43
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +10.06%  "

# Row 44
$ws.Range("E44").Value = "  -9.48%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.04"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -3.22%  "

# Row 46
$ws.Range("E46").Value = "  -0.83%  "

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.241.48"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -2.48%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.55"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -1.42%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.137"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.90%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.53"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +27.39%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.46"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +57.75%  "
